# "Add files via upload" - append two newly-submitted retailer feedback
# rows to the Sheet1 feedback log, and widen the DateSubmitted column's
# number format so it can show the time portion used by the new rows'
# submission timestamps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing DateSubmitted values (E2:E11) were date-only; update their
# display format to include hours/minutes/seconds.
$ws.Range("E2:E11").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# New feedback row from retailer fdf12335-2c41-40b8-a607-920ff9af1019
# about order add6407f-8225-46ae-be70-e5a3c9a9b5c7.
$ws.Range("A12").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("B12").Value = "add6407f-8225-46ae-be70-e5a3c9a9b5c7"
$ws.Range("C12").Value = "General Feedback"
$ws.Range("D12").Value = "hghghglkkkoj"
$ws.Range("E12").Value = "2025-08-07 14:05:55"

# Second new feedback row, same retailer/order, different feedback type.
$ws.Range("A13").Value = "fdf12335-2c41-40b8-a607-920ff9af1019"
$ws.Range("B13").Value = "add6407f-8225-46ae-be70-e5a3c9a9b5c7"
$ws.Range("C13").Value = "Product Issue"
$ws.Range("D13").Value = "fdvdfv"
$ws.Range("E13").Value = "2025-08-07 23:13:44"
